# Updated cryptos list on Tue May 30 08:46:54 UTC 2023 with GitHub Actions
#
# Note: several "Price" values look numeric (e.g. "313.08") but must stay
# plain text, matching the workbook's inline-string cells. A leading
# apostrophe forces Excel to keep the value as text instead of silently
# coercing it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.845.62"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.905.77"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'313.08"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "'0.5057"
$ws.Range("E7").Value = "  +5.15%  "
$ws.Range("D8").Value = "'0.3813"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").Value = "'0.07273"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").Value = "'0.9069"
$ws.Range("E10").Value = "  -2.73%  "
$ws.Range("D11").Value = "'20.91"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.941.75"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07651"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "'5.489"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").Value = "'91.85"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").Value = "'0.000008714"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "27.875.63"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").Value = "'5.157"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "'6.569"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("D24").Value = "'153.67"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("D25").Value = "'1.877"
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("D26").Value = "'2.224"
$ws.Range("E26").Value = "  +4.56%  "
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").Value = "'115.36"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("D29").Value = "'4.904"
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("D30").Value = "'0.09022"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").Value = "'3.211"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("D32").Value = "'1.221"
$ws.Range("E32").Value = "  -3.05%  "
$ws.Range("D33").Value = "'4.696"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").Value = "'0.7661"
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("D35").Value = "'0.02066"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").Value = "'2.512"
$ws.Range("E36").Value = "  -4.94%  "
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("D38").Value = "'0.5529"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").Value = "'3.012"
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").Value = "'0.05260"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").Value = "'6.869"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").Value = "'8.445"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("D44").Value = "'111.10"
$ws.Range("E44").Value = "  +2.85%  "
$ws.Range("D45").Value = "'10.58"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").Value = "'0.4803"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").Value = "'1.628"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "'67.38"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("D50").Value = "'0.06061"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").Value = "'0.9025"
$ws.Range("E51").Value = "  +0.30%  "
